# Rebrand the "Analytics" and "Monetization" section title slides to
# include the "Unity Gaming Services | " prefix.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shape = $s.Shapes.Item($i)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "Analytics") {
                $tr.Text = "Unity Gaming Services | Analytics"
            } elseif ($tr.Text -eq "Monetization") {
                $tr.Text = "Unity Gaming Services | Monetization"
            }
        }
    }
}
